# updated Planning KW 12
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill in the new "Curr. Est.", "Effort" and "Remain" columns for the two task rows
$ws.Range("C2").Value = "1h"
$ws.Range("D2").Value = "2h"
$ws.Range("E2").Value = "2h"

$ws.Range("C3").Value = "2h"
$ws.Range("D3").Value = "2h"
$ws.Range("E3").Value = "4h"

# Update the selected/active cell to E4
$ws.Range("E4").Select()
